$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1395.4445
$ws.Range("I2").Value = 937.6923
$ws.Range("K2").Value = 937.6923
$ws.Range("M2").Value = -824.6923
$ws.Range("H21").Value = 16950
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("H23").Value = 16950
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("H28").Value = 67776.664
$ws.Range("I28").Value = 91768.82000000001
$ws.Range("J28").Value = 1798.25
$ws.Range("K28").Value = 91768.82000000001
$ws.Range("L28").Value = 1798.25
$ws.Range("M28").Value = -91283.82000000001
$ws.Range("N28").Value = -2768.25
$ws.Range("H51").Value = 10736.223
$ws.Range("I51").Value = 15999.333
$ws.Range("K51").Value = 15999.333
$ws.Range("M51").Value = -15515.333
$ws.Range("H74").Value = 15239.9
$ws.Range("I74").Value = 12066.5
$ws.Range("K74").Value = 12066.5
$ws.Range("M74").Value = -11130.5
$ws.Range("H77").Value = 15239.9
$ws.Range("I77").Value = 12066.5
$ws.Range("K77").Value = 60332.5
$ws.Range("M77").Value = -55652.5
$ws.Range("H100").Value = 1349.4286
$ws.Range("I100").Value = 709.4
$ws.Range("K100").Value = 709.4
$ws.Range("M100").Value = -168.4
$ws.Range("H116").Value = 3561.6667
$ws.Range("I116").Value = 2068
$ws.Range("K116").Value = 2068
$ws.Range("M116").Value = 1374
$ws.Range("H138").Value = 3108.1562
$ws.Range("I138").Value = 2214.6667
$ws.Range("J138").Value = 4256.9287
$ws.Range("K138").Value = 6644.000100000001
$ws.Range("L138").Value = 12770.7861
$ws.Range("M138").Value = -1504.000100000001
$ws.Range("N138").Value = -23050.7861
$ws.Range("M21").ClearContents()
$ws.Range("M23").ClearContents()

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4386.5386
$ws.Range("I74").Value = 2111.3635
$ws.Range("K74").Value = 2111.3635
$ws.Range("M74").Value = -1237.3635
$ws.Range("H77").Value = 4386.5386
$ws.Range("I77").Value = 2111.3635
$ws.Range("K77").Value = 10556.8175
$ws.Range("M77").Value = -6188.817499999999
$ws.Range("H104").Value = 28332.666
$ws.Range("J104").Value = 28332.666
$ws.Range("L104").Value = 28332.666
$ws.Range("N104").Value = -35320.666
$ws.Range("H122").Value = 5530.857
$ws.Range("I122").Value = 5368.7
$ws.Range("J122").Value = 5936.25
$ws.Range("K122").Value = 16106.1
$ws.Range("L122").Value = 17808.75
$ws.Range("M122").Value = -13656.1
$ws.Range("N122").Value = -22708.75
$ws.Range("H123").Value = 68111.11
$ws.Range("J123").Value = 68111.11
$ws.Range("L123").Value = 68111.11
$ws.Range("N123").Value = -77911.11

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 690.63336
$ws.Range("I94").Value = 668.2
$ws.Range("K94").Value = 668.2
$ws.Range("M94").Value = -217.2
$ws.Range("H99").Value = 1549.25
$ws.Range("I99").Value = 1519.4
$ws.Range("K99").Value = 1519.4
$ws.Range("M99").Value = -21.40000000000009
$ws.Range("H107").Value = 2359
$ws.Range("I107").Value = 1723.625
$ws.Range("K107").Value = 1723.625
$ws.Range("M107").Value = 196.375
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 89698.7
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 89698.7
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 89698.7
$ws.Range("N59").Value = -91988.7
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("M131").ClearContents()
$ws.Range("N131").ClearContents()

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 44.88889
$ws.Range("J12").Value = 71.8
$ws.Range("L12").Value = 215.4
$ws.Range("N12").Value = -561.4
$ws.Range("H37").Value = 139994
$ws.Range("J37").Value = 139994
$ws.Range("L37").Value = 419982
$ws.Range("N37").Value = -420206
$ws.Range("H92").Value = 93
$ws.Range("J92").Value = 93
$ws.Range("L92").Value = 279
$ws.Range("N92").Value = -2775
$ws.Range("H94").Value = 10154
$ws.Range("J94").Value = 13585.333
$ws.Range("L94").Value = 40755.999
$ws.Range("N94").Value = -42107.999
$ws.Range("H137").Value = 2881.1765
$ws.Range("I137").Value = 1650.2
$ws.Range("J137").Value = 4639.7144
$ws.Range("K137").Value = 4950.6
$ws.Range("L137").Value = 13919.1432
$ws.Range("M137").Value = 149.3999999999996
$ws.Range("N137").Value = -24119.1432
$ws.Range("H140").Value = 3780.1372
$ws.Range("I140").Value = 6874.1665
$ws.Range("J140").Value = 2092.4849
$ws.Range("K140").Value = 20622.4995
$ws.Range("L140").Value = 6277.4547
$ws.Range("M140").Value = -15442.4995
$ws.Range("N140").Value = -16637.4547
$ws.Range("H141").Value = 7487.5
$ws.Range("I141").Value = 7487.5
$ws.Range("K141").Value = 22462.5
$ws.Range("M141").Value = -17282.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 49662.668
$ws.Range("J92").Value = 49662.668
$ws.Range("L92").Value = 49662.668
$ws.Range("N92").Value = -53406.668
$ws.Range("H122").Value = 7301.9165
$ws.Range("I122").Value = 5262.3
$ws.Range("J122").Value = 17500
$ws.Range("K122").Value = 15786.9
$ws.Range("L122").Value = 52500
$ws.Range("M122").Value = -13336.9
$ws.Range("N122").Value = -57400
$ws.Range("H123").Value = 55747.4
$ws.Range("J123").Value = 55747.4
$ws.Range("L123").Value = 55747.4
$ws.Range("N123").Value = -60647.4
$ws.Range("H126").Value = 4223.1665
$ws.Range("I126").Value = 4747.1113
$ws.Range("J126").Value = 2651.3333
$ws.Range("K126").Value = 14241.3339
$ws.Range("L126").Value = 7953.999899999999
$ws.Range("M126").Value = -11771.3339
$ws.Range("N126").Value = -12893.9999

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7198.375
$ws.Range("I7").Value = 6764.6665
$ws.Range("K7").Value = 6764.6665
$ws.Range("M7").Value = -6652.6665
$ws.Range("H55").Value = 1234.9412
$ws.Range("I55").Value = 165.2
$ws.Range("K55").Value = 165.2
$ws.Range("M55").Value = 7.800000000000011
$ws.Range("H126").Value = 7198.375
$ws.Range("I126").Value = 6764.6665
$ws.Range("K126").Value = 20293.9995
$ws.Range("M126").Value = -17823.9995
$ws.Range("H136").Value = 7746.027
$ws.Range("J136").Value = 10377.777
$ws.Range("L136").Value = 31133.331
$ws.Range("N136").Value = -36233.331

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 1000
$ws.Range("J18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("N18").Value = -1346
$ws.Range("H96").Value = 4192
$ws.Range("J96").Value = 4997
$ws.Range("L96").Value = 4997
$ws.Range("N96").Value = -7743
$ws.Range("H113").Value = 1084.4166
$ws.Range("I113").Value = 1334.75
$ws.Range("J113").Value = 583.75
$ws.Range("K113").Value = 4004.25
$ws.Range("L113").Value = 1751.25
$ws.Range("M113").Value = -1834.25
$ws.Range("N113").Value = -6091.25
$ws.Range("H122").Value = 4226.5474
$ws.Range("I122").Value = 4252.2886
$ws.Range("J122").Value = 2888
$ws.Range("K122").Value = 12756.8658
$ws.Range("L122").Value = 8664
$ws.Range("M122").Value = -10306.8658
$ws.Range("N122").Value = -13564
$ws.Range("H132").Value = 3338.8823
$ws.Range("I132").Value = 2422.5625
$ws.Range("K132").Value = 7267.6875
$ws.Range("M132").Value = -4737.6875
